## Apply "more work towards final product" edit.
## This adds a handful of new data cells to the worksheet:
##  - a "D" column value (carrier) for the four practice-pair rows (2-5)
##  - a "J" column value (pair_kind = unique_video/unique_audio) for the
##    four generic-pair rows (6-9)
##  - new "C"/"D" (kind/carrier) values for rows 14-21, which previously
##    only held the running stimulus number in column A

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows: carrier (column D) mirrors the pair_kind already stored
# in column K for these rows.
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic rows 6-9: add the new pair_kind (column J) values.
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Rows 14-21: populate kind (C) and carrier (D) for the remaining
# stimulus numbers (9-16).
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
